$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.86%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.10%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.496"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.45%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08026"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.37%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.091"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "10.26%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9531"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.85%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1149"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.25%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1914"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.66%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "18.50%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09911"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.33%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04883"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12.47%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1065"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.10%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001275"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.81%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04080"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005979"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.08%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.373"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.78%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.403"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.39%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.99%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1380"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.12%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.30%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001274"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.03%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004359"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.60%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.80%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003744"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.17%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02590"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.54%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05794"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "6.30%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007561"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.34%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1404"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.42%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007314"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.52%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.08%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009064"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.91%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007007"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.25%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.03%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.17%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003530"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "55.45%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003537"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.45%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.03%"
